# "Hard works - add students"
# Extends the attendance sheet with two additional class-day columns (G,H),
# bumps each student's missed-class counter, and flips the green/red
# pattern-fill styles (style index 2 <-> style index 3) so the coloring
# stays consistent with the newly inserted days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Stash copies of the three cell formats we need as re-usable
#    "donors" in scratch cells well outside the used range, BEFORE any
#    of the source cells are touched.
# ---------------------------------------------------------------------
$donorGreen = $ws.Range("Z1")   # style index 2 (the fill that starts life as green)
$donorRed   = $ws.Range("Z2")   # style index 3 (the fill that starts life as red)
$donorDate  = $ws.Range("Z3")   # style index 1 (date format used in row 1)

$ws.Range("B2").Copy()
$donorGreen.PasteSpecial($xlPasteFormats)

$ws.Range("D2").Copy()
$donorRed.PasteSpecial($xlPasteFormats)

$ws.Range("C1").Copy()
$donorDate.PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 2) New header cells for the two additional class days (row 1).
# ---------------------------------------------------------------------
$ws.Range("G1").Value = 45637
$ws.Range("H1").Value = 45637

$donorDate.Copy()
$ws.Range("G1:H1").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 3) Regular student rows (2-8, 10-14): two more tracked class days,
#    one more missed class, and the existing C/D/E/F columns swap
#    their fill between green (style 2) and red (style 3).
# ---------------------------------------------------------------------
$regularRows = @(2,3,4,5,6,7,8,10,11,12,13,14)

foreach ($r in $regularRows) {
    $ws.Range("B$r").Value = 2

    $donorRed.Copy()
    $ws.Range("C$r").PasteSpecial($xlPasteFormats)

    $donorGreen.Copy()
    $ws.Range("D$r").PasteSpecial($xlPasteFormats)

    $donorRed.Copy()
    $ws.Range("E$r`:G$r").PasteSpecial($xlPasteFormats)

    $donorGreen.Copy()
    $ws.Range("H$r").PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# 4) Row 9 (Морозов) only ever tracked two class days before this edit;
#    it now gains two tracked days (E,F) instead of (G,H).
# ---------------------------------------------------------------------
$ws.Range("B9").Value = 1

$donorRed.Copy()
$ws.Range("B9").PasteSpecial($xlPasteFormats)
$ws.Range("C9:E9").PasteSpecial($xlPasteFormats)

$donorGreen.Copy()
$ws.Range("F9").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 5) Totals row.
# ---------------------------------------------------------------------
$ws.Range("A17").Value = 6

# ---------------------------------------------------------------------
# 6) Remove the scratch donor cells so they don't leak into the sheet.
# ---------------------------------------------------------------------
$donorGreen.Clear()
$donorRed.Clear()
$donorDate.Clear()
